# Export not included label with blank string by default
#
# 1. Rename the misspelled defined name "dis_value_values" -> "disp_value_values"
# 2. Add a new "Gender Not Used" lookup entry (index 0) to the Gender cut on
#    the Lookups sheet, and shift the existing "male"/"female" entries for
#    the question_code / result_type mapping columns down by one so that the
#    "not used" placeholder keeps index 0 while real answers now start at 1.

$wb = $excel.ActiveWorkbook

# --- 1. Fix the typo'd defined name -------------------------------------------------
$definedName = $wb.Names.Item("dis_value_values")
$definedName.Name = "disp_value_values"

# --- 2. Update the Lookups sheet table ----------------------------------------------
$lookups = $wb.Worksheets.Item("Lookups")

# Gender cut (column F/G): insert the "Gender Not Used" placeholder as the
# first (index 0) entry, push "male"/"female" to indexes 1/2.
$lookups.Range("F2").Value = "Gender Not Used"
$lookups.Range("G2").Value = 0

$lookups.Range("F3").Value = "male"
$lookups.Range("G3").Value = 1

$lookups.Range("F4").Value = "female"
$lookups.Range("G4").Value = 2

# question_code cut (column H/I): keep male/female but start counting at 1
# now that index 0 is reserved for "not used".
$lookups.Range("H2").Value = "male"
$lookups.Range("I2").Value = 1

$lookups.Range("H3").Value = "female"
$lookups.Range("I3").Value = 2

# result_type cut (column J/K): same shift as question_code.
$lookups.Range("J2").Value = "male"
$lookups.Range("K2").Value = 1

$lookups.Range("J3").Value = "female"
$lookups.Range("K3").Value = 2
